# Add a new "21-jul" date column (AI) to the daily-tracking table, appending
# the next day's values after the existing "20-jul" column (AH).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("AI1").Value = "21-jul"

# New column's data values, one per data row (2-18).
$aiValues = @{
    2  = 0
    3  = 8.7689015631202185
    4  = 11.521334511659333
    5  = 24.179073320231613
    6  = 0
    7  = 14.265769208504937
    8  = 8.9100214653227088
    9  = 18.18247688837014
    10 = 21.113639532588557
    11 = 14.257185852491116
    12 = 0
    13 = 13.288712091839649
    14 = 0
    15 = 0
    16 = 10.772324494682191
    17 = 0
    18 = 0
}

foreach ($row in $aiValues.Keys) {
    $ws.Cells.Item($row, 35).Value = $aiValues[$row]
}

# The hidden "helper" column block (previously K:S) grows to K:AB as more
# date columns pile up after the visible table.
$ws.Range("T1:AB1").EntireColumn.Hidden = $true
$ws.Range("T1:AB1").EntireColumn.ColumnWidth = 0

# The sheet was left scrolled one column over with AK8 as the active cell.
$ws.Range("AK8").Select()
